$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 changes
$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 2.8
$ws.Range("I3").Value = 4.2
$ws.Range("M3").Value = 1.17
$ws.Range("N3").Value = 5
$ws.Range("O3").Value = 1.67
$ws.Range("P3").Value = 2.1
$ws.Range("Y3").Value = 11
$ws.Range("AC3").Value = 5
$ws.Range("AE3").Value = 23
$ws.Range("AN3").Value = 3.75
$ws.Range("AR3").Value = 101
$ws.Range("AV3").Value = 101

# Row 6 changes
$ws.Range("G6").Value = 2.25
$ws.Range("I6").Value = 3.25
$ws.Range("J6").Value = 3
$ws.Range("L6").Value = 4
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.62
$ws.Range("S6").Value = 1.5
$ws.Range("T6").Value = 2.5
$ws.Range("U6").Value = 1.91
$ws.Range("V6").Value = 1.8
$ws.Range("W6").Value = 6.5
$ws.Range("X6").Value = 10
$ws.Range("Z6").Value = 21
$ws.Range("AA6").Value = 21
$ws.Range("AC6").Value = 8
$ws.Range("AE6").Value = 17
$ws.Range("AF6").Value = 51
$ws.Range("AH6").Value = 15
$ws.Range("AI6").Value = 12
$ws.Range("AJ6").Value = 34
$ws.Range("AK6").Value = 29
$ws.Range("AM6").Value = 800
$ws.Range("AN6").Value = 4.33
$ws.Range("AO6").Value = 13
$ws.Range("AT6").Value = 2.5
$ws.Range("AU6").Value = 8.5
$ws.Range("AW6").Value = 5
$ws.Range("AX6").Value = 19
$ws.Range("AY6").Value = 29
$ws.Range("AZ6").Value = 67
$ws.Range("BB6").Value = 251

# Row 7 changes
$ws.Range("G7").Value = 1.48
$ws.Range("H7").Value = 4.2
$ws.Range("I7").Value = 7
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("U7").Value = 2.1
$ws.Range("V7").Value = 1.67
$ws.Range("X7").Value = 6.5
$ws.Range("Y7").Value = 8.5
$ws.Range("Z7").Value = 9.5
$ws.Range("AB7").Value = 29
$ws.Range("AC7").Value = 10
$ws.Range("AH7").Value = 34
$ws.Range("AO7").Value = 7.5
$ws.Range("AW7").Value = 8
$ws.Range("BA7").Value = 151
$ws.Range("BB7").Value = 351
